# Add a new employee record (row 10) to the employee list, mirroring the
# existing rows' layout (columns B:H), then leave the selection where the
# user would land after typing the row (one row below, column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "0010"
$ws.Range("C10").Value = "Paresh"
$ws.Range("D10").Value = "Sonaawane"
$ws.Range("E10").Value = "Test Analyst"
$ws.Range("F10").Value = "Full-Time Permanent"
$ws.Range("G10").Value = "Dev"
$ws.Range("H10").Value = "Sunil Dolwani"

[void]$ws.Range("G11").Select()
